$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38: politeness_score was stored as text "3"; normalize it to a real
# number to match the rest of column B.
$ws.Range("B38").Value = 3

# New row 39 - additional annotation for Ying Tang.
$ws.Range("A39").Value = "Ying Tang"

# politeness_score for this row was entered as text (keeps the leading
# quote so it stays a string "2" instead of being coerced to a number).
$ws.Range("B39").Value = "'2"

$ws.Range("C39").Value = "simplistic"
$ws.Range("D39").Value = "CRT"
$ws.Range("E39").Value = "MET"
$ws.Range("F39").Value = "0c8a854c-e7df-48dd-93a0-b6771319a745"
$ws.Range("G39").Value = "H1Ww66x0-_annotated.xlsx"
$ws.Range("H39").Value = "- the proposed approach to maintain the budget is simplistic"
